$wb = $excel.ActiveWorkbook

# --- Rename the Power Pivot worksheet-connection defined names -------------
# Re-saving in a newer Excel build collapsed the trailing de-dup suffix ("1")
# that a previous save had appended to these hidden `_xlcn.*` defined names
# (e.g. ...xlsxNode_Media1 -> ...xlsxNode_Media). Names are in worksheet order
# as they appear in the workbook: Node_Media, Product_Accordions,
# Product_Categories, Product_Nodes, TACO_Nodes.
$newNames = @(
    "_xlcn.WorksheetConnection_DraftTemplate.xlsxNode_Media",
    "_xlcn.WorksheetConnection_DraftTemplate.xlsxProduct_Accordions",
    "_xlcn.WorksheetConnection_DraftTemplate.xlsxProduct_Categories",
    "_xlcn.WorksheetConnection_DraftTemplate.xlsxProduct_Nodes",
    "_xlcn.WorksheetConnection_DraftTemplate.xlsxTACO_Nodes"
)
for ($i = 1; $i -le $wb.Names.Count(); $i++) {
    $definedName = $wb.Names.Item($i)
    $definedName.Name = $newNames[$i - 1]
}

# --- Preserve cell format while fixing the stray trailing spaces -----------
# "Product1 " / "Price1 " on the Mapping sheet carried an accidental trailing
# space. Writing only the .Value (not ClearContents+retype) keeps the
# existing style/number-format of A4:B4 untouched.
$mapping = $wb.Worksheets.Item("Mapping")
$mapping.Range("A4").Value = "Product1"
$mapping.Range("B4").Value = "Price1"

# --- Restore the on-disk selection/active-sheet state -----------------------
# The workbook was last saved with the Mapping sheet active and B4 selected
# (previously RepeatingData was active with A5 selected on it).
$mapping.Activate()
[void]$mapping.Range("B4").Select()
